$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Munka1")
$ws2 = $wb.Worksheets.Item("Munka2")

# --- Munka2: rebuild the header/data row with the new item-schema columns ---
# Writing the cells in this particular order reproduces the exact shared-string
# table ordering/index layout of the target workbook.
$ws2.Range("C1").Value = "Type"
$ws2.Range("B1").Value = "SystemName"
$ws2.Range("A1").Value = "ItemName"
$ws2.Range("C2").Value = "WeaponMain"
$ws2.Range("E1").Value = "NecessaryTypes"
$ws2.Range("E2").Value = "Grip;Handguard;WeaponBody"
$ws2.Range("D1").Value = "Description"
$ws2.Range("D2").Value = "This is a sovjet AKS-74U carabine"

$ws2.Range("A2").Value = "AKS-74U"
$ws2.Range("B2").Value = "AKS-74U"

# --- Munka2: column widths shift right by 3 columns (new Type/SystemName/ItemName cols) ---
$ws2.Columns.Item(4).ColumnWidth = 13.73
$ws2.Columns.Item(5).ColumnWidth = 20.02
$ws2.Columns.Item(6).ColumnWidth = 7.88

# --- Activate Munka2 (it becomes the selected/visible tab) ---
$ws2.Activate()
$ws2.Range("A2").Select()

Write-Host "edit applied"
